$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "Bigquery procedure 를 이용하여 recursion 함수 만들기"
$ws.Range("E32").Value = "https://dodonam.tistory.com/317"

$ws.Range("D46").Value = "[LG전자] 2021년 04월, 생물정보학(Bioinformatics 채용), 인체/미생물 유전자 및 바이오인포매틱스 전문가 모집"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/391"

$ws.Range("D50").Value = "Overleaf를 이용한 논문작성"
$ws.Range("E50").Value = "http://incredible.egloos.com/7515316"

$ws.Range("D51").Value = "[sqlite3] 윈도우 10에 sqlite3 설치 및 환경변수 path 설정하기"
$ws.Range("E51").Value = "https://bskyvision.com/1174"
